# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Hoja1: refresh the "Conversión del día" note with the new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$currentText = $ws1.Range("A1").Value()
$updatedText = $currentText.Replace(
    "✅ 1000 Bs = 9.77 = 41039.16 pesos`n✅ 41039.16 pesos = 9.7 = 953.67 Bs",
    "✅ 1000 Bs = 9.9 = 41741.09 pesos`n✅ 41741.09 pesos = 9.9 = 969.81 Bs"
)
$ws1.Range("A1").Value = $updatedText

# --- tasas: refresh the Binance/transfi rate lookup cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 101
$ws2.Range("O10").Value = 4215.85
$ws2.Range("N12").Value = 4217.97
$ws2.Range("O12").Value = 98
